# Auto-generated COM-interop script applying the BRVM recommandations update
$wb = $excel.ActiveWorkbook
$wsReco = $wb.Worksheets.Item("Recommandations")
$wsYtd = $wb.Worksheets.Item("Top_YTD")

# --- Sheet "Recommandations": update rows 2-44 in place (values/labels shifted) ---
$wsReco.Range("A2").Value = "BRVM - SERVICES PUBLICS"
$wsReco.Range("B2").Value = 0
$wsReco.Range("C2").Value = 8
$wsReco.Range("D2").Value = 3328.81
$wsReco.Range("E2").Value = 104.29
$wsReco.Range("F2").Value = "🟡 Observer"
$wsReco.Range("G2").Value = "➖ Neutre"

$wsReco.Range("A3").Value = "AIR LIQUIDE CI"
$wsReco.Range("B3").Value = 0
$wsReco.Range("C3").Value = 4
$wsReco.Range("D3").Value = 2760
$wsReco.Range("E3").Value = 690
$wsReco.Range("F3").Value = "🟡 Observer"
$wsReco.Range("G3").Value = "➖ Neutre"

$wsReco.Range("A4").Value = "NEI-CEDA CI"
$wsReco.Range("B4").Value = 0
$wsReco.Range("C4").Value = 4
$wsReco.Range("D4").Value = 2675
$wsReco.Range("E4").Value = 680
$wsReco.Range("F4").Value = "🟡 Observer"
$wsReco.Range("G4").Value = "➖ Neutre"

$wsReco.Range("A5").Value = "BRVM - AUTRES SECTEURS"
$wsReco.Range("B5").Value = 0
$wsReco.Range("C5").Value = 4
$wsReco.Range("D5").Value = 2403.32
$wsReco.Range("E5").Value = 566.83
$wsReco.Range("F5").Value = "🟡 Observer"
$wsReco.Range("G5").Value = "➖ Neutre"

$wsReco.Range("A6").Value = "BRVM - DISTRIBUTION"
$wsReco.Range("B6").Value = 0
$wsReco.Range("C6").Value = 4
$wsReco.Range("D6").Value = 1707.05
$wsReco.Range("E6").Value = 429.88
$wsReco.Range("F6").Value = "🟡 Observer"
$wsReco.Range("G6").Value = "➖ Neutre"

$wsReco.Range("A7").Value = "BRVM - TRANSPORT"
$wsReco.Range("B7").Value = 0
$wsReco.Range("C7").Value = 4
$wsReco.Range("D7").Value = 1450.27
$wsReco.Range("E7").Value = 358.59
$wsReco.Range("F7").Value = "🟡 Observer"
$wsReco.Range("G7").Value = "➖ Neutre"

$wsReco.Range("A8").Value = "BRVM - AGRICULTURE"
$wsReco.Range("B8").Value = 0
$wsReco.Range("C8").Value = 4
$wsReco.Range("D8").Value = 1445.72
$wsReco.Range("E8").Value = 360.73
$wsReco.Range("F8").Value = "🟡 Observer"
$wsReco.Range("G8").Value = "➖ Neutre"

$wsReco.Range("A9").Value = "BRVM - CONSOMMATION DISCRETIONNAIRE"
$wsReco.Range("B9").Value = 0
$wsReco.Range("C9").Value = 4
$wsReco.Range("D9").Value = 549.83
$wsReco.Range("E9").Value = 136.28
$wsReco.Range("F9").Value = "🟡 Observer"
$wsReco.Range("G9").Value = "➖ Neutre"

$wsReco.Range("A10").Value = "BRVM-PRESTIGE"
$wsReco.Range("B10").Value = 0
$wsReco.Range("C10").Value = 4
$wsReco.Range("D10").Value = 547.04
$wsReco.Range("E10").Value = 137.27
$wsReco.Range("F10").Value = "🟡 Observer"
$wsReco.Range("G10").Value = "➖ Neutre"

$wsReco.Range("A11").Value = "BRVM - FINANCES"
$wsReco.Range("B11").Value = 0
$wsReco.Range("C11").Value = 4
$wsReco.Range("D11").Value = 542.75
$wsReco.Range("E11").Value = 135.92
$wsReco.Range("F11").Value = "🟡 Observer"
$wsReco.Range("G11").Value = "➖ Neutre"

$wsReco.Range("A12").Value = "BRVM - SERVICES FINANCIERS"
$wsReco.Range("B12").Value = 0
$wsReco.Range("C12").Value = 4
$wsReco.Range("D12").Value = 533.4
$wsReco.Range("E12").Value = 133.58
$wsReco.Range("F12").Value = "🟡 Observer"
$wsReco.Range("G12").Value = "➖ Neutre"

$wsReco.Range("A13").Value = "BRVM - INDUSTRIELS"
$wsReco.Range("B13").Value = 0
$wsReco.Range("C13").Value = 4
$wsReco.Range("D13").Value = 481.68
$wsReco.Range("E13").Value = 122.54
$wsReco.Range("F13").Value = "🟡 Observer"
$wsReco.Range("G13").Value = "➖ Neutre"

$wsReco.Range("A14").Value = "BRVM - ENERGIE"
$wsReco.Range("B14").Value = 0
$wsReco.Range("C14").Value = 4
$wsReco.Range("D14").Value = 430.42
$wsReco.Range("E14").Value = 107.73
$wsReco.Range("F14").Value = "🟡 Observer"
$wsReco.Range("G14").Value = "➖ Neutre"

$wsReco.Range("A15").Value = "BRVM - INDUSTRIE                  (**)"
$wsReco.Range("B15").Value = 0
$wsReco.Range("C15").Value = 2
$wsReco.Range("D15").Value = 414.76
$wsReco.Range("E15").Value = 207.68
$wsReco.Range("F15").Value = "🟡 Observer"
$wsReco.Range("G15").Value = "➖ Neutre"

$wsReco.Range("A16").Value = "BRVM-PRINCIPAL                    (**)"
$wsReco.Range("B16").Value = 0
$wsReco.Range("C16").Value = 2
$wsReco.Range("D16").Value = 387.23
$wsReco.Range("E16").Value = 193.83
$wsReco.Range("F16").Value = "🟡 Observer"
$wsReco.Range("G16").Value = "➖ Neutre"

$wsReco.Range("A17").Value = "BRVM - TELECOMMUNICATIONS"
$wsReco.Range("B17").Value = 0
$wsReco.Range("C17").Value = 4
$wsReco.Range("D17").Value = 379.79
$wsReco.Range("E17").Value = 95.56
$wsReco.Range("F17").Value = "🟡 Observer"
$wsReco.Range("G17").Value = "➖ Neutre"

$wsReco.Range("A18").Value = "BRVM - CONSOMMATION DE BASE         (**)"
$wsReco.Range("B18").Value = 0
$wsReco.Range("C18").Value = 2
$wsReco.Range("D18").Value = 372.63
$wsReco.Range("E18").Value = 187.07
$wsReco.Range("F18").Value = "🟡 Observer"
$wsReco.Range("G18").Value = "➖ Neutre"

$wsReco.Range("A19").Value = "BRVM - INDUSTRIE                (**)"
$wsReco.Range("B19").Value = 0
$wsReco.Range("C19").Value = 1
$wsReco.Range("D19").Value = 206.88
$wsReco.Range("E19").Value = 206.88
$wsReco.Range("F19").Value = "🟡 Observer"
$wsReco.Range("G19").Value = "➖ Neutre"

$wsReco.Range("A20").Value = "BRVM-PRINCIPAL                  (**)"
$wsReco.Range("B20").Value = 0
$wsReco.Range("C20").Value = 1
$wsReco.Range("D20").Value = 193.58
$wsReco.Range("E20").Value = 193.58
$wsReco.Range("F20").Value = "🟡 Observer"
$wsReco.Range("G20").Value = "➖ Neutre"

$wsReco.Range("A21").Value = "BRVM - CONSOMMATION DE BASE              (**)"
$wsReco.Range("B21").Value = 0
$wsReco.Range("C21").Value = 1
$wsReco.Range("D21").Value = 187.55
$wsReco.Range("E21").Value = 187.55
$wsReco.Range("F21").Value = "🟡 Observer"
$wsReco.Range("G21").Value = "➖ Neutre"

$wsReco.Range("A22").Value = "FILTISAC CI (FTSC)"
$wsReco.Range("B22").Value = 4
$wsReco.Range("C22").Value = 0
$wsReco.Range("D22").Value = 29.49
$wsReco.Range("E22").Value = 7.34
$wsReco.Range("F22").Value = "🟢 Achat"
$wsReco.Range("G22").Value = "✅ Renforcer"

$wsReco.Range("A23").Value = "ORAGROUP TOGO (ORGT)"
$wsReco.Range("B23").Value = 2
$wsReco.Range("C23").Value = 1
$wsReco.Range("D23").Value = 11.39
$wsReco.Range("E23").Value = 7.41
$wsReco.Range("F23").Value = "🟡 Observer"
$wsReco.Range("G23").Value = "👀 À surveiller"

$wsReco.Range("A24").Value = "BANK OF AFRICA BF (BOABF)"
$wsReco.Range("B24").Value = 1
$wsReco.Range("C24").Value = 0
$wsReco.Range("D24").Value = 4.85
$wsReco.Range("E24").Value = 4.85
$wsReco.Range("F24").Value = "🟡 Observer"
$wsReco.Range("G24").Value = "➖ Neutre"

$wsReco.Range("A25").Value = "SODE CI (SDCC)"
$wsReco.Range("B25").Value = 1
$wsReco.Range("C25").Value = 0
$wsReco.Range("D25").Value = 4.35
$wsReco.Range("E25").Value = 4.35
$wsReco.Range("F25").Value = "🟡 Observer"
$wsReco.Range("G25").Value = "➖ Neutre"

$wsReco.Range("A26").Value = "VIVO ENERGY CI (SHEC)"
$wsReco.Range("B26").Value = 2
$wsReco.Range("C26").Value = 0
$wsReco.Range("D26").Value = 4.19
$wsReco.Range("E26").Value = 2.92
$wsReco.Range("F26").Value = "🟡 Observer"
$wsReco.Range("G26").Value = "➖ Neutre"

$wsReco.Range("A27").Value = "SETAO CI (STAC)"
$wsReco.Range("B27").Value = 1
$wsReco.Range("C27").Value = 1
$wsReco.Range("D27").Value = 3.52
$wsReco.Range("E27").Value = -1.9
$wsReco.Range("F27").Value = "🟡 Observer"
$wsReco.Range("G27").Value = "👀 À surveiller"

$wsReco.Range("A28").Value = "ONATEL BF (ONTBF)"
$wsReco.Range("B28").Value = 1
$wsReco.Range("C28").Value = 0
$wsReco.Range("D28").Value = 2.61
$wsReco.Range("E28").Value = 2.61
$wsReco.Range("F28").Value = "🟡 Observer"
$wsReco.Range("G28").Value = "➖ Neutre"

$wsReco.Range("A29").Value = "NSIA BANQUE COTE D'IVOIRE (NSBC)"
$wsReco.Range("B29").Value = 1
$wsReco.Range("C29").Value = 0
$wsReco.Range("D29").Value = 2.5
$wsReco.Range("E29").Value = 2.5
$wsReco.Range("F29").Value = "🟡 Observer"
$wsReco.Range("G29").Value = "➖ Neutre"

$wsReco.Range("A30").Value = "BICI CI (BICC)"
$wsReco.Range("B30").Value = 1
$wsReco.Range("C30").Value = 0
$wsReco.Range("D30").Value = 2.2
$wsReco.Range("E30").Value = 2.2
$wsReco.Range("F30").Value = "🟡 Observer"
$wsReco.Range("G30").Value = "➖ Neutre"

$wsReco.Range("A31").Value = "CFAO MOTORS CI (CFAC)"
$wsReco.Range("B31").Value = 1
$wsReco.Range("C31").Value = 1
$wsReco.Range("D31").Value = 1.39
$wsReco.Range("E31").Value = 3.17
$wsReco.Range("F31").Value = "🟡 Observer"
$wsReco.Range("G31").Value = "👀 À surveiller"

$wsReco.Range("A32").Value = "SICABLE CI (CABC)"
$wsReco.Range("B32").Value = 1
$wsReco.Range("C32").Value = 1
$wsReco.Range("D32").Value = 0.25
$wsReco.Range("E32").Value = 5.2
$wsReco.Range("F32").Value = "🟡 Observer"
$wsReco.Range("G32").Value = "👀 À surveiller"

$wsReco.Range("A33").Value = "TOTAL"
$wsReco.Range("B33").Value = 0
$wsReco.Range("C33").Value = 4
$wsReco.Range("D33").Value = 0
$wsReco.Range("E33").Value = 0
$wsReco.Range("F33").Value = "🟡 Observer"
$wsReco.Range("G33").Value = "➖ Neutre"

$wsReco.Range("A34").Value = "SAPH CI (SPHC)"
$wsReco.Range("B34").Value = 1
$wsReco.Range("C34").Value = 1
$wsReco.Range("D34").Value = -0.32
$wsReco.Range("E34").Value = 5.94
$wsReco.Range("F34").Value = "🟡 Observer"
$wsReco.Range("G34").Value = "👀 À surveiller"

$wsReco.Range("A35").Value = "BERNABE CI (BNBC)"
$wsReco.Range("B35").Value = 1
$wsReco.Range("C35").Value = 1
$wsReco.Range("D35").Value = -1.01
$wsReco.Range("E35").Value = 6.45
$wsReco.Range("F35").Value = "🟡 Observer"
$wsReco.Range("G35").Value = "👀 À surveiller"

$wsReco.Range("A36").Value = "AFRICA GLOBAL LOGISTICS CI (SDSC)"
$wsReco.Range("B36").Value = 0
$wsReco.Range("C36").Value = 1
$wsReco.Range("D36").Value = -1.68
$wsReco.Range("E36").Value = -1.68
$wsReco.Range("F36").Value = "🟡 Observer"
$wsReco.Range("G36").Value = "➖ Neutre"

$wsReco.Range("A37").Value = "TOTALENERGIES MARKETING CI (TTLC)"
$wsReco.Range("B37").Value = 0
$wsReco.Range("C37").Value = 1
$wsReco.Range("D37").Value = -2.04
$wsReco.Range("E37").Value = -2.04
$wsReco.Range("F37").Value = "🟡 Observer"
$wsReco.Range("G37").Value = "➖ Neutre"

$wsReco.Range("A38").Value = "UNILEVER CI (UNLC)"
$wsReco.Range("B38").Value = 1
$wsReco.Range("C38").Value = 1
$wsReco.Range("D38").Value = -2.24
$wsReco.Range("E38").Value = -7.5
$wsReco.Range("F38").Value = "🟡 Observer"
$wsReco.Range("G38").Value = "👀 À surveiller"

$wsReco.Range("A39").Value = "NEI-CEDA CI (NEIC)"
$wsReco.Range("B39").Value = 0
$wsReco.Range("C39").Value = 1
$wsReco.Range("D39").Value = -2.94
$wsReco.Range("E39").Value = -2.94
$wsReco.Range("F39").Value = "🟡 Observer"
$wsReco.Range("G39").Value = "➖ Neutre"

$wsReco.Range("A40").Value = "TRACTAFRIC MOTORS CI (PRSC)"
$wsReco.Range("B40").Value = 1
$wsReco.Range("C40").Value = 2
$wsReco.Range("D40").Value = -6.57
$wsReco.Range("E40").Value = 7.46
$wsReco.Range("F40").Value = "🟡 Observer"
$wsReco.Range("G40").Value = "👀 À surveiller"

$wsReco.Range("A41").Value = "SICOR CI (SICC)"
$wsReco.Range("B41").Value = 0
$wsReco.Range("C41").Value = 1
$wsReco.Range("D41").Value = -6.91
$wsReco.Range("E41").Value = -6.91
$wsReco.Range("F41").Value = "🟡 Observer"
$wsReco.Range("G41").Value = "➖ Neutre"

$wsReco.Range("A42").Value = "UNIWAX CI (UNXC)"
$wsReco.Range("B42").Value = 0
$wsReco.Range("C42").Value = 2
$wsReco.Range("D42").Value = -10.79
$wsReco.Range("E42").Value = -7.48
$wsReco.Range("F42").Value = "🟡 Observer"
$wsReco.Range("G42").Value = "➖ Neutre"

$wsReco.Range("A43").Value = "SUCRIVOIRE (SCRC)"
$wsReco.Range("B43").Value = 0
$wsReco.Range("C43").Value = 3
$wsReco.Range("D43").Value = -13.69
$wsReco.Range("E43").Value = -3.43
$wsReco.Range("F43").Value = "🔴 Vente"
$wsReco.Range("G43").Value = "⚠️ Risque de décrochage"

$wsReco.Range("A44").Value = "LOTERIE NATIONALE DU BENIN (LNBB)"
$wsReco.Range("B44").Value = 0
$wsReco.Range("C44").Value = 2
$wsReco.Range("D44").Value = -14.92
$wsReco.Range("E44").Value = -7.45
$wsReco.Range("F44").Value = "🟡 Observer"
$wsReco.Range("G44").Value = "➖ Neutre"

# Remove now-obsolete trailing rows 45-47 (table shrank from 47 to 44 data rows)
$wsReco.Rows.Item(45).Delete()
$wsReco.Rows.Item(45).Delete()
$wsReco.Rows.Item(45).Delete()

# --- Sheet "Top_YTD": update Progression YTD (%) values ---
$wsYtd.Range("B2").Value = 8204097.68
$wsYtd.Range("B3").Value = 389338.4
$wsYtd.Range("B4").Value = 348999.92
$wsYtd.Range("B5").Value = 239948.85
$wsYtd.Range("B6").Value = 76889.67
$wsYtd.Range("B7").Value = 45679.91
$wsYtd.Range("B8").Value = 45224.4
$wsYtd.Range("B9").Value = 3078.54
$wsYtd.Range("B10").Value = 3042.18
$wsYtd.Range("B11").Value = 2985.64

Write-Output "Update complete"
